# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# file "be4ba6f3-03a4-44a7-a837-338a32d96c05.md" has now been handed back
# (previously it was "Ready for handoff"). As part of regenerating the
# report, the rows in each sheet are re-sorted alphabetically by source
# file name, which moves the be4ba6f3 file from the last row to the first
# data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws1.Range("B2").Value = "e2e\be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws1.Range("C2").Value = ".md"
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("G2").Value = "2016-09-06 07:15:46"

$ws1.Range("A3").Value = "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md"
$ws1.Range("B3").Value = "e2e\ffff778366a4-ef15-40b7-a639-bcc0e5053014.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-09-06 07:14:00"

$ws1.Range("A4").Value = "ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md"
$ws1.Range("B4").Value = "e2e\ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md"
$ws1.Range("C4").Value = ".md"
$ws1.Range("E4").Value = "Handed back: in sync with en-US"
$ws1.Range("F4").Value = "Handed back: in sync with en-US"
$ws1.Range("G4").Value = "2016-09-06 07:14:00"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/be4ba6f3-03a4-44a7-a837-338a32d96c05.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\be4ba6f3-03a4-44a7-a837-338a32d96c05.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd520aa1a543a0628de92e6165f6f122a0efa04/e2e/ffff778366a4-ef15-40b7-a639-bcc0e5053014.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\ffff778366a4-ef15-40b7-a639-bcc0e5053014.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("F2").Value = "False"
$ws2.Range("G2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.3c1fd758bf93f305cb73dbbd541807dbc06c39f3.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-06 07:15:40"
$ws2.Range("I2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws2.Range("J2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.3c1fd758bf93f305cb73dbbd541807dbc06c39f3.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-06 07:16:20"

$ws2.Range("A3").Value = "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("F3").Value = "False"
$ws2.Range("G3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-06 07:13:55"
$ws2.Range("I3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md"
$ws2.Range("J3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-06 07:14:20"

$ws2.Range("A4").Value = "ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md"
$ws2.Range("C4").Value = "Handed back: in sync with en-US"
$ws2.Range("F4").Value = "True"
$ws2.Range("G4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-09-06 07:13:55"
$ws2.Range("I4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md"
$ws2.Range("J4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-09-06 07:14:20"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/be4ba6f3-03a4-44a7-a837-338a32d96c05.md", [System.Type]::Missing, [System.Type]::Missing, "be4ba6f3-03a4-44a7-a837-338a32d96c05.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/be4ba6f3-03a4-44a7-a837-338a32d96c05.md", [System.Type]::Missing, [System.Type]::Missing, "be4ba6f3-03a4-44a7-a837-338a32d96c05.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd520aa1a543a0628de92e6165f6f122a0efa04/e2e/ffff778366a4-ef15-40b7-a639-bcc0e5053014.md", [System.Type]::Missing, [System.Type]::Missing, "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/189e3b7d7c59d010c2d3a8b82f0284650592d78b/e2e/a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md", [System.Type]::Missing, [System.Type]::Missing, "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md", [System.Type]::Missing, [System.Type]::Missing, "ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/189e3b7d7c59d010c2d3a8b82f0284650592d78b/e2e/a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md", [System.Type]::Missing, [System.Type]::Missing, "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("F2").Value = "False"
$ws3.Range("G2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.3c1fd758bf93f305cb73dbbd541807dbc06c39f3.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-06 07:15:46"
$ws3.Range("I2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
$ws3.Range("J2").Value = "be4ba6f3-03a4-44a7-a837-338a32d96c05.3c1fd758bf93f305cb73dbbd541807dbc06c39f3.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-06 07:16:27"

$ws3.Range("A3").Value = "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("F3").Value = "False"
$ws3.Range("G3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-06 07:14:00"
$ws3.Range("I3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md"
$ws3.Range("J3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-06 07:14:28"

$ws3.Range("A4").Value = "ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md"
$ws3.Range("C4").Value = "Handed back: in sync with en-US"
$ws3.Range("F4").Value = "True"
$ws3.Range("G4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$ws3.Range("H4").Value = "2016-09-06 07:14:00"
$ws3.Range("I4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md"
$ws3.Range("J4").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$ws3.Range("K4").Value = "2016-09-06 07:14:28"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/be4ba6f3-03a4-44a7-a837-338a32d96c05.md", [System.Type]::Missing, [System.Type]::Missing, "be4ba6f3-03a4-44a7-a837-338a32d96c05.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/be4ba6f3-03a4-44a7-a837-338a32d96c05.md", [System.Type]::Missing, [System.Type]::Missing, "be4ba6f3-03a4-44a7-a837-338a32d96c05.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd520aa1a543a0628de92e6165f6f122a0efa04/e2e/ffff778366a4-ef15-40b7-a639-bcc0e5053014.md", [System.Type]::Missing, [System.Type]::Missing, "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/844af212fa29bc071bef4864a5c7d51b09d3dbdf/e2e/a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md", [System.Type]::Missing, [System.Type]::Missing, "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b07a4dcf887791d995ef1a6d11ba17b79af1fb9/e2e/ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md", [System.Type]::Missing, [System.Type]::Missing, "ffffff0a6e813f-e7e8-4587-91f5-ee1e7a63c19c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/844af212fa29bc071bef4864a5c7d51b09d3dbdf/e2e/a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md", [System.Type]::Missing, [System.Type]::Missing, "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md") | Out-Null

Write-Host "Report regenerated for handback of be4ba6f3-03a4-44a7-a837-338a32d96c05.md"
